# Auto-generated edit script applying the cryptos.xlsx price/volume/coin update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $Val)
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "37.414.81"
$ws.Range("E2").Value = "  +6.19%  "
# Row 3
Set-TextValue "D3" "2.039.04"
$ws.Range("E3").Value = "  +8.13%  "
# Row 4
$ws.Range("E4").Value = "  +0.02%  "
# Row 5
Set-TextValue "D5" "253.65"
$ws.Range("E5").Value = "  +3.58%  "
# Row 6
Set-TextValue "D6" "0.694"
$ws.Range("E6").Value = "  +0.84%  "
# Row 7
$ws.Range("E7").Value = "  +0.03%  "
# Row 8
Set-TextValue "D8" "46.73"
$ws.Range("E8").Value = "  +9.62%  "
# Row 9
$ws.Range("E9").Value = "  +8.54%  "
# Row 10
Set-TextValue "D10" "57.43"
$ws.Range("E10").Value = "  +4.63%  "
# Row 11
Set-TextValue "D11" "0.0765"
$ws.Range("E11").Value = "  +3.45%  "
# Row 12
$ws.Range("E12").Value = "  +2.53%  "
# Row 13
Set-TextValue "D13" "15.47"
$ws.Range("E13").Value = "  +12.77%  "
# Row 14
Set-TextValue "D14" "0.840"
# Row 15
Set-TextValue "D15" "2.330.64"
$ws.Range("E15").Value = "  +7.96%  "
# Row 16
Set-TextValue "D16" "5.21"
$ws.Range("E16").Value = "  +5.42%  "
# Row 17
Set-TextValue "D17" "2.032.00"
$ws.Range("E17").Value = "  +7.75%  "
# Row 18
Set-TextValue "D18" "37.439.12"
$ws.Range("E18").Value = "  +6.27%  "
# Row 19
Set-TextValue "D19" "75.40"
$ws.Range("E19").Value = "  +3.20%  "
# Row 20
Set-TextValue "D20" "0.0₃0863"
$ws.Range("E20").Value = "  +5.01%  "
# Row 21
Set-TextValue "D21" "13.76"
$ws.Range("E21").Value = "  +7.94%  "
# Row 22
Set-TextValue "D22" "254.04"
$ws.Range("E22").Value = "  +4.35%  "
# Row 23
$ws.Range("E23").Value = "  +2.75%  "
# Row 24
Set-TextValue "D24" "0.999"
$ws.Range("E24").Value = "  -0.13%  "
# Row 25
$ws.Range("E25").Value = "  -4.71%  "
# Row 26
Set-TextValue "D26" "170.60"
$ws.Range("E26").Value = "  +1.98%  "
# Row 27
$ws.Range("E27").Value = "  -0.47%  "
# Row 28
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D28" "20.71"
$ws.Range("E28").Value = "  +13.74%  "
# Row 29
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D29" "8.96"
$ws.Range("E29").Value = "  +5.92%  "
# Row 30
Set-TextValue "D30" "0.129"
$ws.Range("E30").Value = "  +2.43%  "
# Row 31
Set-TextValue "D31" "22.88"
$ws.Range("E31").Value = "  +73.99%  "
# Row 32
$ws.Range("E32").Value = "  +7.35%  "
# Row 33
Set-TextValue "D33" "0.0617"
$ws.Range("E33").Value = "  +4.37%  "
# Row 34
$ws.Range("E34").Value = "  +4.65%  "
# Row 35
Set-TextValue "D35" "0.0896"
$ws.Range("E35").Value = "  +25.34%  "
# Row 36
$ws.Range("E36").Value = "  +0.06%  "
# Row 37
Set-TextValue "D37" "1.89"
$ws.Range("E37").Value = "  +1.04%  "
# Row 38
Set-TextValue "D38" "2.30"
$ws.Range("E38").Value = "  +19.00%  "
# Row 39
Set-TextValue "D39" "0.902"
$ws.Range("E39").Value = "  +6.74%  "
# Row 40
$ws.Range("E40").Value = "  -0.16%  "
# Row 41
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D41" "0.0229"
$ws.Range("E41").Value = "  +4.10%  "
# Row 42
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D42" "103.00"
$ws.Range("E42").Value = "  +5.35%  "
# Row 43
Set-TextValue "D43" "17.40"
$ws.Range("E43").Value = "  +1.97%  "
# Row 44
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D44" "2.92"
$ws.Range("E44").Value = "  +21.03%  "
# Row 45
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D45" "1.14"
$ws.Range("E45").Value = "  +7.04%  "
# Row 46
Set-TextValue "D46" "1.373.27"
$ws.Range("E46").Value = "  +3.76%  "
# Row 47
$ws.Range("E47").Value = "  +4.86%  "
# Row 48
$ws.Range("E48").Value = "  +2.90%  "
# Row 49
Set-TextValue "D49" "2.86"
$ws.Range("E49").Value = "  +4.63%  "
# Row 50
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D50" "2.220.68"
$ws.Range("E50").Value = "  +7.84%  "
# Row 51
Set-TextValue "D51" "6.64"
$ws.Range("E51").Value = "  +6.42%  "
